# Generate Report for Handoff
# Updates the localization-status report: the handoff batch id moves from
# 9c16adeb-27c1-4662-baa6-517a55aa3656 to 4d419704-5596-4684-b5f2-bb781cf8d6f1
# (new .md handoff file, new .xlf target files) and refreshes the handoff
# timestamps for the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldId = "9c16adeb-27c1-4662-baa6-517a55aa3656"
$newId = "4d419704-5596-4684-b5f2-bb781cf8d6f1"
$oldHash = "c01d58be974233ac5c323f5d483d7e21dc6fc146"
$newHash = "3536d8ea407a5a49725ae9d35ab119d8a61a723f"

$newMdName = "$newId.md"
$newZhXlf = "$newId.$newHash.zh-cn.xlf"
$newDeXlf = "$newId.$newHash.de-de.xlf"

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/e16ff43ec89e004a1cb3c04c9d1832c6ae261ecb/e2e/$oldId.md"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/865d01904cbe32e78161fe7f6e0daab8e1396fb0/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$oldId.$oldHash.zh-cn.xlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/673a88cf437538d3ebf22df1bc83659048e3e0e7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$oldId.$oldHash.de-de.xlf"

function Update-HyperlinkDisplay {
    # Positional params: (worksheet, cellRef, target, newText) — this runtime's
    # PowerShell host does not bind named (-Foo bar) arguments on user
    # functions, so keep call sites positional.
    param($Worksheet, $CellRef, $Target, $NewText)

    # Re-point every hyperlink anchored on this cell so its stored `display`
    # text matches the refreshed cell value (the link target itself does not
    # change). Rebuilding the whole collection is what actually rewrites the
    # existing <hyperlink> entry in place instead of appending a stale one.
    $range = $Worksheet.Range($CellRef)
    $range.Value2 = $NewText

    $toDelete = @()
    foreach ($hl in $Worksheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $range.Address()) {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }

    $Worksheet.Hyperlinks.Add($range, $Target, [Type]::Missing, [Type]::Missing, $NewText) | Out-Null
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
Update-HyperlinkDisplay $wsOverview "A2" $mdTarget $newMdName
$wsOverview.Range("D2").Value2 = "2016-03-22 05:40:53"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HyperlinkDisplay $wsZhCn "A2" $mdTarget $newMdName
Update-HyperlinkDisplay $wsZhCn "D2" $zhXlfTarget $newZhXlf
$wsZhCn.Range("E2").Value2 = "2016-03-22 05:40:45"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HyperlinkDisplay $wsDeDe "A2" $mdTarget $newMdName
Update-HyperlinkDisplay $wsDeDe "D2" $deXlfTarget $newDeXlf
$wsDeDe.Range("E2").Value2 = "2016-03-22 05:40:53"
